$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 113
$ws_ALC.Range("H113").Value = 20412110
$ws_ALC.Range("I113").Value = 111113370
$ws_ALC.Range("J113").Value = 4327.2
$ws_ALC.Range("K113").Value = 111113370
$ws_ALC.Range("L113").Value = 4327.2
$ws_ALC.Range("M113").Value = -111110116
$ws_ALC.Range("N113").Value = -10835.2

# Row 129
$ws_ALC.Range("H129").Value = 854.7536
$ws_ALC.Range("J129").Value = 877.8280999999999
$ws_ALC.Range("L129").Value = 2633.4843
$ws_ALC.Range("N129").Value = -12633.4843

# Row 137
$ws_ALC.Range("H137").Value = 35092.098
$ws_ALC.Range("J137").Value = 64747.062
$ws_ALC.Range("L137").Value = 194241.186
$ws_ALC.Range("N137").Value = -199341.186

# Row 138
$ws_ALC.Range("H138").Value = 2696.392
$ws_ALC.Range("J138").Value = 3055.4878
$ws_ALC.Range("L138").Value = 9166.463400000001
$ws_ALC.Range("N138").Value = -19446.4634

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 32
$ws_ARM.Range("H32").Value = 18268.188
$ws_ARM.Range("I32").Value = 21924.92
$ws_ARM.Range("K32").Value = 21924.92
$ws_ARM.Range("M32").Value = -21637.92

# Row 45
$ws_ARM.Range("H45").Value = 2599.3333
$ws_ARM.Range("I45").Value = 1965.2354
$ws_ARM.Range("K45").Value = 1965.2354
$ws_ARM.Range("M45").Value = -1588.2354

# Row 61
$ws_ARM.Range("H61").Value = 412221.3
$ws_ARM.Range("I61").Value = 623280.1
$ws_ARM.Range("J61").Value = 4174.2666
$ws_ARM.Range("K61").Value = 623280.1
$ws_ARM.Range("L61").Value = 4174.2666
$ws_ARM.Range("M61").Value = -623068.1
$ws_ARM.Range("N61").Value = -4598.2666

# Row 88
$ws_ARM.Range("H88").Value = 69665.664
$ws_ARM.Range("I88").Value = 1574.3334
$ws_ARM.Range("J88").Value = 115059.89
$ws_ARM.Range("K88").Value = 1574.3334
$ws_ARM.Range("L88").Value = 115059.89
$ws_ARM.Range("M88").Value = -1168.3334
$ws_ARM.Range("N88").Value = -115871.89

# Row 91
$ws_ARM.Range("H91").Value = 69665.664
$ws_ARM.Range("I91").Value = 1574.3334
$ws_ARM.Range("J91").Value = 115059.89
$ws_ARM.Range("K91").Value = 1574.3334
$ws_ARM.Range("L91").Value = 115059.89
$ws_ARM.Range("M91").Value = -170.3334
$ws_ARM.Range("N91").Value = -117867.89

# Row 132
$ws_ARM.Range("H132").Value = 13077.645
$ws_ARM.Range("I132").Value = 1956.3077
$ws_ARM.Range("K132").Value = 5868.9231
$ws_ARM.Range("M132").Value = -3338.9231

# Row 136
$ws_ARM.Range("H136").Value = 412221.3
$ws_ARM.Range("I136").Value = 623280.1
$ws_ARM.Range("J136").Value = 4174.2666
$ws_ARM.Range("K136").Value = 1869840.3
$ws_ARM.Range("L136").Value = 12522.7998
$ws_ARM.Range("M136").Value = -1867290.3
$ws_ARM.Range("N136").Value = -17622.7998

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 86
$ws_BSM.Range("H86").Value = 1719.95
$ws_BSM.Range("I86").Value = 1562.5
$ws_BSM.Range("J86").Value = 2349.75
$ws_BSM.Range("K86").Value = 1562.5
$ws_BSM.Range("L86").Value = 2349.75
$ws_BSM.Range("M86").Value = -439.5
$ws_BSM.Range("N86").Value = -4595.75

# Row 89
$ws_BSM.Range("H89").Value = 1719.95
$ws_BSM.Range("I89").Value = 1562.5
$ws_BSM.Range("J89").Value = 2349.75
$ws_BSM.Range("K89").Value = 7812.5
$ws_BSM.Range("L89").Value = 11748.75
$ws_BSM.Range("M89").Value = -2196.5
$ws_BSM.Range("N89").Value = -22980.75

# Row 99
$ws_BSM.Range("H99").Value = 1422.4166
$ws_BSM.Range("I99").Value = 898.625
$ws_BSM.Range("J99").Value = 2470
$ws_BSM.Range("K99").Value = 898.625
$ws_BSM.Range("L99").Value = 2470
$ws_BSM.Range("M99").Value = 599.375
$ws_BSM.Range("N99").Value = -5466

# Row 134
$ws_BSM.Range("H134").Value = 56551.105
$ws_BSM.Range("I134").Value = 56551.105
$ws_BSM.Range("J134").Value = 0
$ws_BSM.Range("K134").Value = 169653.315
$ws_BSM.Range("L134").Value = 0
$ws_BSM.Range("M134").Value = -167118.315
$ws_BSM.Range("N134").ClearContents()

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 31
$ws_CRP.Range("H31").Value = 9423.386
$ws_CRP.Range("I31").Value = 14559.3545
$ws_CRP.Range("J31").Value = 3299.7307
$ws_CRP.Range("K31").Value = 14559.3545
$ws_CRP.Range("L31").Value = 3299.7307
$ws_CRP.Range("M31").Value = -14264.3545
$ws_CRP.Range("N31").Value = -3889.7307

# Row 34
$ws_CRP.Range("H34").Value = 9423.386
$ws_CRP.Range("I34").Value = 14559.3545
$ws_CRP.Range("J34").Value = 3299.7307
$ws_CRP.Range("K34").Value = 14559.3545
$ws_CRP.Range("L34").Value = 3299.7307
$ws_CRP.Range("M34").Value = -14357.3545
$ws_CRP.Range("N34").Value = -3703.7307

# Row 62
$ws_CRP.Range("H62").Value = 6168.5
$ws_CRP.Range("I62").Value = 4499.5
$ws_CRP.Range("K62").Value = 4499.5
$ws_CRP.Range("M62").Value = -3875.5

# Row 65
$ws_CRP.Range("H65").Value = 6168.5
$ws_CRP.Range("I65").Value = 4499.5
$ws_CRP.Range("K65").Value = 22497.5
$ws_CRP.Range("M65").Value = -19377.5

# Row 112
$ws_CRP.Range("H112").Value = 45700
$ws_CRP.Range("J112").Value = 45700
$ws_CRP.Range("L112").Value = 45700
$ws_CRP.Range("N112").Value = -48654

# Row 122
$ws_CRP.Range("H122").Value = 3500.25
$ws_CRP.Range("I122").Value = 4333.6665
$ws_CRP.Range("J122").Value = 1000
$ws_CRP.Range("K122").Value = 13000.9995
$ws_CRP.Range("L122").Value = 3000
$ws_CRP.Range("M122").Value = -10550.9995
$ws_CRP.Range("N122").Value = -7900

# Row 134
$ws_CRP.Range("H134").Value = 6076.263
$ws_CRP.Range("I134").Value = 732.2941
$ws_CRP.Range("K134").Value = 2196.8823
$ws_CRP.Range("M134").Value = 338.1177000000002

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 38
$ws_CUL.Range("H38").Value = 162
$ws_CUL.Range("I38").Value = 183
$ws_CUL.Range("J38").Value = 134
$ws_CUL.Range("K38").Value = 549
$ws_CUL.Range("L38").Value = 402
$ws_CUL.Range("M38").Value = -202
$ws_CUL.Range("N38").Value = -1096

# Row 68
$ws_CUL.Range("H68").Value = 1286.1875
$ws_CUL.Range("J68").Value = 1279.2903
$ws_CUL.Range("L68").Value = 3837.8709
$ws_CUL.Range("N68").Value = -5459.8709

# Row 71
$ws_CUL.Range("H71").Value = 1286.1875
$ws_CUL.Range("J71").Value = 1279.2903
$ws_CUL.Range("L71").Value = 11513.6127
$ws_CUL.Range("N71").Value = -19625.6127

# Row 131
$ws_CUL.Range("H131").Value = 803.72
$ws_CUL.Range("J131").Value = 817.1042
$ws_CUL.Range("L131").Value = 2451.3126
$ws_CUL.Range("N131").Value = -12531.3126

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 70
$ws_GSM.Range("H70").Value = 15628250
$ws_GSM.Range("I70").Value = 4000
$ws_GSM.Range("J70").Value = 20836334
$ws_GSM.Range("K70").Value = 4000
$ws_GSM.Range("L70").Value = 20836334
$ws_GSM.Range("M70").Value = -3730
$ws_GSM.Range("N70").Value = -20836874

# Row 73
$ws_GSM.Range("H73").Value = 15628250
$ws_GSM.Range("I73").Value = 4000
$ws_GSM.Range("J73").Value = 20836334
$ws_GSM.Range("K73").Value = 4000
$ws_GSM.Range("L73").Value = 20836334
$ws_GSM.Range("M73").Value = -3064
$ws_GSM.Range("N73").Value = -20838206

# Row 122
$ws_GSM.Range("H122").Value = 6999.8
$ws_GSM.Range("I122").Value = 5000
$ws_GSM.Range("J122").Value = 7499.75
$ws_GSM.Range("K122").Value = 15000
$ws_GSM.Range("L122").Value = 22499.25
$ws_GSM.Range("N122").Value = -27399.25
$ws_GSM.Range("M122").Value = -12550

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 16
$ws_LTW.Range("H16").Value = 2667.2856
$ws_LTW.Range("I16").Value = 2667.2856
$ws_LTW.Range("K16").Value = 2667.2856
$ws_LTW.Range("M16").Value = -2497.2856

# Row 22
$ws_LTW.Range("H22").Value = 1641.5
$ws_LTW.Range("J22").Value = 882
$ws_LTW.Range("L22").Value = 882
$ws_LTW.Range("N22").Value = -1472

# Row 27
$ws_LTW.Range("H27").Value = 1641.5
$ws_LTW.Range("J27").Value = 882
$ws_LTW.Range("L27").Value = 882
$ws_LTW.Range("N27").Value = -1096
